$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newFitness = @(
11555, 11555, 11555, 11288, 11288, 11036, 9677, 9677, 9677, 8942, 8942, 8942, 8942, 8942, 8942, 8942, 8934, 8897, 8897, 8897, 8897, 8897, 8897, 8896, 8896, 8896, 8896, 8896, 8896, 8896, 8673, 8673, 8673, 8218, 8218, 8218, 8218, 8218, 8218, 8218, 8218, 8218, 8218, 8218, 8218, 8218, 7845, 7845, 7845, 7845, 7845, 7845, 7845, 7845, 7845, 7845, 7845, 7845, 7845, 7845, 7845, 7845, 7845, 7845, 7845, 7845, 7320, 7320, 7320, 7320, 7320, 7320, 7318, 7318, 7318, 7318, 7318, 7318, 7318, 7318, 7318, 7318, 7318
)

for ($i = 0; $i -lt $newFitness.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 3).Value = $newFitness[$i]
}
